$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting rows 6-9 down to 7-10.
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new market entry.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C6").Value = 'Los Lagos'
$ws.Range("D6").Value = 44523
$ws.Range("D6").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 'Fruta'
$ws.Range("G6").Value = 100103
$ws.Range("H6").Value = 'Frutos de hueso (carozo)'
$ws.Range("I6").Value = 100103003
$ws.Range("J6").Value = 'Damasco'
$ws.Range("K6").Value = 'Castle Brite'
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 500
$ws.Range("N6").Value = 28000
$ws.Range("O6").Value = 28500
$ws.Range("P6").Value = 28250
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("R6").Value = 'Provincia de Limarí'
$ws.Range("S6").Value = 1569
$ws.Range("T6").Value = 18
